$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column B - widen it and switch the "Story" column cells to Arial font
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 71.4140625

$storyRange = $ws.Range("B2:B9")
$storyRange.Font.Name = "Arial"

$ws.Range("B2").Value = "As an unauthorized user I want create a new account so I can buy and sell items."
$ws.Range("B3").Value = "As an unauthorized user I want to login my account so that I can access my info."
$ws.Range("B4").Value = "As an authorized user I want to logout of my account so that a stranger cannot use my account"
$ws.Range("B5").Value = "As an unauthorized user I want to use my google account so that I can use one of my own existing accounts."
$ws.Range("B6").Value = "As an unauthorized user I want to use my Facebook account so that I can use one of my own existing accounts."
$ws.Range("B7").Value = "As a user I want to be able to add an item to the list"
$ws.Range("B8").Value = "As a seller, I want to be able to sell my items in my local area so that relevant buyers can find my items.  "
$ws.Range("B9").Value = "As a new student, I want to be able to navigate over to the book section so that I can buy textbooks for my classes."

# ---------------------------------------------------------------------------
# Column C - Story Priority
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 2
$ws.Range("C7").Value = 2
$ws.Range("C8").Value = 2
$ws.Range("C9").Value = 2

# ---------------------------------------------------------------------------
# Column D
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 4
$ws.Range("D6").Value = 5
$ws.Range("D7").Value = 6
$ws.Range("D8").Value = 7
$ws.Range("D9").Value = 6

# ---------------------------------------------------------------------------
# Column E - Story Status, styled cells ("IP" = In Progress / Neutral style,
# "W" = Waiting, using the 20% - Accent3 style with a thick double border)
# ---------------------------------------------------------------------------
$ipRange = $ws.Range("E2:E5")
$ipRange.Style = "Neutral"
$ws.Range("E2").Value = "IP"
$ws.Range("E3").Value = "IP"
$ws.Range("E4").Value = "IP"
$ws.Range("E5").Value = "IP"

$wRange = $ws.Range("E6:E10")
$wRange.Style = "20% - Accent3"
$wRange.Borders.LineStyle = -4119
$wRange.Borders.Color = 4144959
$ws.Range("E6").Value = "W"
$ws.Range("E7").Value = "W"
$ws.Range("E8").Value = "W"
$ws.Range("E9").Value = "W"
$ws.Range("E10").Value = "W"

# ---------------------------------------------------------------------------
# Column F - Story Points
# ---------------------------------------------------------------------------
$ws.Range("F2").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 2
$ws.Range("F6").Value = 2

# ---------------------------------------------------------------------------
# Row heights - thick double borders on E6:E10 make Excel bump the row
# height of the affected rows (5 through 11) when it marks them with
# thickTop/thickBot.
# ---------------------------------------------------------------------------
$ws.Rows(5).RowHeight = 16
$ws.Rows(6).RowHeight = 16.5
$ws.Rows(7).RowHeight = 16.5
$ws.Rows(8).RowHeight = 16.5
$ws.Rows(9).RowHeight = 16.5
$ws.Rows(10).RowHeight = 16.5
$ws.Rows(11).RowHeight = 16

# ---------------------------------------------------------------------------
# Selection - the author's last selected cell before saving
# ---------------------------------------------------------------------------
$ws.Range("F7").Select()
